{"js": "// Each entry is [oldText, newText] - applied as an exact-match replace\n// against the document body. All values in this document are unique,\n// so a simple search+replace per pair is unambiguous.\nconst replacements = [\n  [\"2024-01-22 Monday\", \"2024-01-23 Tuesday\"],\n  [\"605\u00d73=\", \"909\u00d72=\"],\n  [\"962\u00d75=\", \"885\u00d73=\"],\n  [\"952\u00d79=\", \"791\u00d78=\"],\n  [\"945\u00d72=\", \"621\u00d78=\"],\n  [\"762\u00d79=\", \"964\u00d73=\"],\n  [\"977\u00d78=\", \"118\u00d77=\"],\n  [\"310\u00d78=\", \"730\u00d78=\"],\n  [\"611\u00d78=\", \"776\u00d75=\"],\n  [\"210\u00d77=\", \"436\u00d74=\"],\n  [\"832\u00d76=\", \"934\u00d79=\"],\n  [\"695\u00d72=\", \"940\u00d77=\"],\n  [\"411\u00d79=\", \"911\u00d74=\"],\n  [\"400\u00d72=\", \"338\u00d78=\"],\n  [\"571\u00d74=\", \"678\u00d77=\"],\n  [\"922\u00d72=\", \"415\u00d76=\"],\n  [\"285\u00d78=\", \"856\u00d72=\"],\n  [\"976\u00d73=\", \"830\u00d77=\"],\n  [\"199\u00d78=\", \"334\u00d76=\"],\n  [\"918\u00d77=\", \"142\u00d74=\"],\n  [\"125\u00d77=\", \"541\u00d73=\"],\n  [\"210\u00d72=\", \"516\u00d73=\"],\n  [\"453\u00d73=\", \"562\u00d72=\"],\n  [\"644\u00d72=\", \"805\u00d77=\"],\n  [\"638\u00d78=\", \"668\u00d78=\"],\n  [\"369\u00d76=\", \"281\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date and every multiplication expression in the table.\n# All \"old\" values are unique in the document, so a Find/Replace pass\n# per pair (ReplaceAll) is unambiguous and only touches the one run.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{old=\"2024-01-22 Monday\"; new=\"2024-01-23 Tuesday\"},\n    @{old=\"605\u00d73=\"; new=\"909\u00d72=\"},\n    @{old=\"962\u00d75=\"; new=\"885\u00d73=\"},\n    @{old=\"952\u00d79=\"; new=\"791\u00d78=\"},\n    @{old=\"945\u00d72=\"; new=\"621\u00d78=\"},\n    @{old=\"762\u00d79=\"; new=\"964\u00d73=\"},\n    @{old=\"977\u00d78=\"; new=\"118\u00d77=\"},\n    @{old=\"310\u00d78=\"; new=\"730\u00d78=\"},\n    @{old=\"611\u00d78=\"; new=\"776\u00d75=\"},\n    @{old=\"210\u00d77=\"; new=\"436\u00d74=\"},\n    @{old=\"832\u00d76=\"; new=\"934\u00d79=\"},\n    @{old=\"695\u00d72=\"; new=\"940\u00d77=\"},\n    @{old=\"411\u00d79=\"; new=\"911\u00d74=\"},\n    @{old=\"400\u00d72=\"; new=\"338\u00d78=\"},\n    @{old=\"571\u00d74=\"; new=\"678\u00d77=\"},\n    @{old=\"922\u00d72=\"; new=\"415\u00d76=\"},\n    @{old=\"285\u00d78=\"; new=\"856\u00d72=\"},\n    @{old=\"976\u00d73=\"; new=\"830\u00d77=\"},\n    @{old=\"199\u00d78=\"; new=\"334\u00d76=\"},\n    @{old=\"918\u00d77=\"; new=\"142\u00d74=\"},\n    @{old=\"125\u00d77=\"; new=\"541\u00d73=\"},\n    @{old=\"210\u00d72=\"; new=\"516\u00d73=\"},\n    @{old=\"453\u00d73=\"; new=\"562\u00d72=\"},\n    @{old=\"644\u00d72=\"; new=\"805\u00d77=\"},\n    @{old=\"638\u00d78=\"; new=\"668\u00d78=\"},\n    @{old=\"369\u00d76=\"; new=\"281\u00d78=\"}\n)\n\n$wdReplaceAll = 2\n\nforeach ($rep in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $rep.old\n    $find.Replacement.Text = $rep.new\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $wdReplaceAll) | Out-Null\n}\n"}
